$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Columns F (Confidence %) and G (Odds) store numeric-looking values as plain
# text in the source feed, so force text formatting before writing, then drop
# back to the Normal style so no stray number-format style lingers on the cells.
$ws.Range("F2:G21").NumberFormat = "@"

$ws.Range("A2").Value = 'Liverpool v Aston Villa'
$ws.Range("B2").Value = 'Liverpool'
$ws.Range("C2").Value = 'England Premier League'
$ws.Range("D2").Value = '2025-11-01T20:00:00.000Z'
$ws.Range("E2").Value = '57/111 Win Tips'
$ws.Range("F2").Value = '51'
$ws.Range("G2").Value = '1.70'

$ws.Range("A3").Value = 'Tottenham v Chelsea'
$ws.Range("B3").Value = 'Chelsea'
$ws.Range("C3").Value = 'England Premier League'
$ws.Range("D3").Value = '2025-11-01T17:30:00.000Z'
$ws.Range("E3").Value = '55/115 Win Tips'
$ws.Range("F3").Value = '48'
$ws.Range("G3").Value = '2.50'

$ws.Range("A4").Value = 'West Ham v Newcastle'
$ws.Range("B4").Value = 'Newcastle'
$ws.Range("C4").Value = 'England Premier League'
$ws.Range("D4").Value = '2025-11-02T14:00:00.000Z'
$ws.Range("E4").Value = '41/54 Win Tips'
$ws.Range("F4").Value = '76'
$ws.Range("G4").Value = '1.67'

$ws.Range("A5").Value = 'Man City v Bournemouth'
$ws.Range("B5").Value = 'Man City'
$ws.Range("C5").Value = 'England Premier League'
$ws.Range("D5").Value = '2025-11-02T16:30:00.000Z'
$ws.Range("E5").Value = '39/52 Win Tips'
$ws.Range("F5").Value = '75'
$ws.Range("G5").Value = '1.53'

$ws.Range("A6").Value = 'Real Madrid v Valencia'
$ws.Range("B6").Value = 'Real Madrid'
$ws.Range("C6").Value = 'Spain Primera Liga'
$ws.Range("D6").Value = '2025-11-01T20:00:00.000Z'
$ws.Range("E6").Value = '27/34 Win Tips'
$ws.Range("F6").Value = '79'
$ws.Range("G6").Value = '1.18'

$ws.Range("A7").Value = 'Sunderland v Everton'
$ws.Range("B7").Value = 'Sunderland'
$ws.Range("C7").Value = 'England Premier League'
$ws.Range("D7").Value = '2025-11-03T20:00:00.000Z'
$ws.Range("E7").Value = '24/35 Win Tips'
$ws.Range("F7").Value = '69'
$ws.Range("G7").Value = '2.88'

$ws.Range("A8").Value = 'Napoli v Como'
$ws.Range("B8").Value = 'Napoli'
$ws.Range("C8").Value = 'Italy Serie A'
$ws.Range("D8").Value = '2025-11-01T17:00:00.000Z'
$ws.Range("E8").Value = '23/29 Win Tips'
$ws.Range("F8").Value = '79'
$ws.Range("G8").Value = '2.00'

$ws.Range("A9").Value = 'Bayern Munich v Bayer Leverkusen'
$ws.Range("B9").Value = 'Bayern Munich'
$ws.Range("C9").Value = 'Germany Bundesliga I'
$ws.Range("D9").Value = '2025-11-01T17:30:00.000Z'
$ws.Range("E9").Value = '22/27 Win Tips'
$ws.Range("F9").Value = '81'
$ws.Range("G9").Value = '1.22'

$ws.Range("A10").Value = 'Feyenoord v FC Volendam'
$ws.Range("B10").Value = 'Feyenoord'
$ws.Range("C10").Value = 'Netherlands Eredivisie'
$ws.Range("D10").Value = '2025-11-01T19:00:00.000Z'
$ws.Range("E10").Value = '16/18 Win Tips'
$ws.Range("F10").Value = '89'
$ws.Range("G10").Value = '1.12'

$ws.Range("A11").Value = 'Monaco v Paris FC'
$ws.Range("B11").Value = 'Monaco'
$ws.Range("C11").Value = 'France Ligue 1'
$ws.Range("D11").Value = '2025-11-01T18:00:00.000Z'
$ws.Range("E11").Value = '15/17 Win Tips'
$ws.Range("F11").Value = '88'
$ws.Range("G11").Value = '1.67'

$ws.Range("A12").Value = 'Cremonese v Juventus'
$ws.Range("B12").Value = 'Juventus'
$ws.Range("C12").Value = 'Italy Serie A'
$ws.Range("D12").Value = '2025-11-01T19:45:00.000Z'
$ws.Range("E12").Value = '14/22 Win Tips'
$ws.Range("F12").Value = '64'
$ws.Range("G12").Value = '1.50'

$ws.Range("A13").Value = 'Club Brugge v Dender'
$ws.Range("B13").Value = 'Club Brugge'
$ws.Range("C13").Value = 'Belgium First Division A'
$ws.Range("D13").Value = '2025-11-01T17:15:00.000Z'
$ws.Range("E13").Value = '12/14 Win Tips'
$ws.Range("F13").Value = '86'
$ws.Range("G13").Value = '1.25'

$ws.Range("A14").Value = 'Guimaraes v Benfica'
$ws.Range("B14").Value = 'Benfica'
$ws.Range("C14").Value = 'Portugal Primeira Liga'
$ws.Range("D14").Value = '2025-11-01T20:30:00.000Z'
$ws.Range("E14").Value = '9/11 Win Tips'
$ws.Range("F14").Value = '82'
$ws.Range("G14").Value = '1.50'

$ws.Range("A15").Value = 'Galatasaray v Trabzonspor'
$ws.Range("B15").Value = 'Galatasaray'
$ws.Range("C15").Value = 'Turkey Super Lig'
$ws.Range("D15").Value = '2025-11-01T17:00:00.000Z'
$ws.Range("E15").Value = '8/10 Win Tips'
$ws.Range("F15").Value = '80'
$ws.Range("G15").Value = '1.44'

$ws.Range("A16").Value = 'Auxerre v Marseille'
$ws.Range("B16").Value = 'Marseille'
$ws.Range("C16").Value = 'France Ligue 1'
$ws.Range("D16").Value = '2025-11-01T20:05:00.000Z'
$ws.Range("E16").Value = '8/11 Win Tips'
$ws.Range("F16").Value = '73'
$ws.Range("G16").Value = '1.70'

$ws.Range("A17").Value = 'Barcelona v Elche'
$ws.Range("B17").Value = 'Barcelona'
$ws.Range("C17").Value = 'Spain Primera Liga'
$ws.Range("D17").Value = '2025-11-02T17:30:00.000Z'
$ws.Range("E17").Value = '8/8 Win Tips'
$ws.Range("F17").Value = '100'
$ws.Range("G17").Value = '1.22'

$ws.Range("A18").Value = 'Virtus Entella v Empoli'
$ws.Range("B18").Value = 'Draw'
$ws.Range("C18").Value = 'Italy Serie B'
$ws.Range("D18").Value = '2025-11-01T16:15:00.000Z'
$ws.Range("E18").Value = '7/10 Win Tips'
$ws.Range("F18").Value = '70'
$ws.Range("G18").Value = '3.15'

$ws.Range("A19").Value = 'FC Copenhagen v FC Fredericia'
$ws.Range("B19").Value = 'FC Copenhagen'
$ws.Range("C19").Value = 'Denmark Superligaen'
$ws.Range("D19").Value = '2025-11-01T17:00:00.000Z'
$ws.Range("E19").Value = '7/10 Win Tips'
$ws.Range("F19").Value = '70'
$ws.Range("G19").Value = '1.18'

$ws.Range("A20").Value = 'Olympiacos v Aris Salonika'
$ws.Range("B20").Value = 'Olympiacos'
$ws.Range("C20").Value = 'Greece Super League 1'
$ws.Range("D20").Value = '2025-11-01T18:00:00.000Z'
$ws.Range("E20").Value = '7/9 Win Tips'
$ws.Range("F20").Value = '78'
$ws.Range("G20").Value = '1.30'

$ws.Range("A21").Value = 'Verona v Inter Milan'
$ws.Range("B21").Value = 'Inter Milan'
$ws.Range("C21").Value = 'Italy Serie A'
$ws.Range("D21").Value = '2025-11-02T11:30:00.000Z'
$ws.Range("E21").Value = '7/7 Win Tips'
$ws.Range("F21").Value = '100'
$ws.Range("G21").Value = '1.44'

# Restore default styling on the numeric-text columns (keeps text type, no format pollution)
$ws.Range("F2:G21").Style = "Normal"
